# Update "想去人数" (want-to-go count) figures that changed between the
# previous data snapshot and the one generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 3215
$wsExhibit.Range("F6").Value = 2058
$wsExhibit.Range("F7").Value = 27
$wsExhibit.Range("F8").Value = 74

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 18

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3215
$wsAll.Range("F3").Value = 18
$wsAll.Range("F7").Value = 2058
$wsAll.Range("F8").Value = 27
$wsAll.Range("F9").Value = 74
